# Daily attendance processing - 2025-10-09 22:50:32
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "System, system, backup@backdoor.com"
$ws.Range("G7").Value = "System, admin@admin.com"
$ws.Range("G10").Value = "System, dnasr281@gmail.com"
$ws.Range("L10").Value = "'70.5%"
$ws.Range("S17").Value = "'61.8%"
$ws.Range("S18").Value = "'76.4%"
$ws.Range("S19").Value = "'76.1%"
$ws.Range("G29").Value = "System, system, backup@backdoor.com"
$ws.Range("G34").Value = "System, admin@admin.com"
$ws.Range("G37").Value = "System, dnasr281@gmail.com"
$ws.Range("G56").Value = "System, system, backup@backdoor.com"
$ws.Range("H59").Value = "31/55"
$ws.Range("G61").Value = "System, admin@admin.com"
$ws.Range("H61").Value = "40/55"
$ws.Range("G64").Value = "System, dnasr281@gmail.com"
$ws.Range("H64").Value = "29/55"
$ws.Range("H85").Value = "45/56"
$ws.Range("H86").Value = "38/56"
$ws.Range("G87").Value = "System, dnasr281@gmail.com"
$ws.Range("H87").Value = "25/56"
$ws.Range("G88").Value = "System, dnasr281@gmail.com"
$ws.Range("H88").Value = "42/56"
$ws.Range("G90").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("H112").Value = "54/55"
$ws.Range("G113").Value = "System, dnasr281@gmail.com"
$ws.Range("H113").Value = "26/55"
$ws.Range("G114").Value = "System, dnasr281@gmail.com"
$ws.Range("G116").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G139").Value = "System, dnasr281@gmail.com"
$ws.Range("G140").Value = "System, dnasr281@gmail.com"
$ws.Range("G142").Value = "admin@admin.com, dnasr281@gmail.com"
